$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content for the "hardware section / sensor section" additions
# (order matters for shared-string table indices, matching the original edit order)
$ws.Range("B2").Value = "Exceptional Geo-referencing tools"
$ws.Range("C2").Value = "Form building tool requires significant learning time; Trial study could not easilly configure a listing of lengths for review in the field"
$ws.Range("C3").Value = "API connectors to external devices requires the Azure license"
$ws.Range("B3").Value = "*Connectors to databases including SQL & OneDrive is No-Code (easy to set up) and functions smoothly;   Able to set up on-the-fly editing, listings and frequency plots to review lengths in the field"

# Row 3 grows considerably taller to fit the new long text
$ws.Rows.Item(3).RowHeight = 156.75

# Apply the wrap-text style to every populated cell (new style index 1)
$ws.Range("A1").WrapText = $true
$ws.Range("B1").WrapText = $true
$ws.Range("C1").WrapText = $true

$ws.Range("A2").WrapText = $true
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true

$ws.Range("A3").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

$ws.Range("A4").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("A6").WrapText = $true
$ws.Range("A7").WrapText = $true

$ws.Range("A8").WrapText = $true
$ws.Range("B8").WrapText = $true

# Selection moves to C3, matching the last-edited cell
$ws.Range("C3").Select()
